$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    ("Play Diamond Queen Free Slot by IGT"). Target structure:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Read our review of ... free spins.</w:t></w:r>
#      </w:p>
# ------------------------------------------------------------------

# Locate the bold "Play Diamond Queen Free Slot by IGT" paragraph near
# the bottom of the document - it already has the exact empty-run +
# bold-run pattern we need, so we copy it as a template (it gets
# deleted later anyway per the diff).
$boldSourceIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Diamond Queen Free Slot by IGT`r") {
        $boldSourceIndex = $i
        break
    }
}
$boldSource = $d.Paragraphs.Item($boldSourceIndex)
$boldSource.Range.Copy()

$titlePara = $d.Paragraphs.First
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
# Set style to Normal BEFORE pasting so the paste logic doesn't treat
# the inherited Heading1 bold as redundant and strip the explicit <w:b/>.
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Append the plain (non-bold) remainder as a separate run.
$metaPara.Range.InsertAfter(": Read our review of Diamond Queen, a magical and elegant online slot by IGT. Play for free and trigger the Mystical Diamond Bonus for extra Wilds and free spins.")

# Replace the pasted "Play Diamond Queen Free Slot by IGT" text (still
# bold) with "Meta description", keeping it inside the same bold run.
$metaRange = $metaPara.Range
$metaRange.Find.Execute("Play Diamond Queen Free Slot by IGT", $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Near the bottom: remove the bold "Play Diamond Queen Free Slot by
#    IGT" paragraph entirely, and rewrite the italic "Read our review..."
#    paragraph's text into the DALL-E prompt (keeping the leading empty
#    run + italic run structure, and without smart-quote mangling).
# ------------------------------------------------------------------

$boldSource = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Diamond Queen Free Slot by IGT`r") {
        $boldSource = $p
        break
    }
}
$boldSource.Range.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$searchRange = $lastPara.Range.Duplicate
$searchRange.Find.Execute("Read our review of Diamond Queen, a magical and elegant online slot by IGT. Play for free and trigger the Mystical Diamond Bonus for extra Wilds and free spins.") | Out-Null
$searchRange.Text = 'Prompt for DALLE: Create a feature image for the online slot game "Diamond Queen" featuring a happy Maya warrior with glasses in a cartoon style. The image should be vibrant and eye-catching, with the Maya warrior holding a large diamond scepter, surrounded by precious jewels and enchanted forest elements. The image should convey a sense of magic and fantasy, while also highlighting the diamond theme of the game. The Maya warrior should be depicted with a joyful expression and be wearing glasses, emphasizing the technological aspect of the game. Overall, the image should appeal to players who are looking for an exciting and magical gaming experience combined with cutting-edge technology.'

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
